# Smartline-IUCNGET crosswalk update
# - Fixes "Soft rock shoes" -> "Soft rock shores" typo (header sheet + SSSOM rows)
# - Adds new SSSOM mapping rows 7-11 (Muddy/Coarse sediment/Undifferentiated sediment/
#   Coral coasts/No stability classification) authored by Rebecca Jordan
# - Updates header sheet's smartline source URL (row 3) and highlights it
# - Updates selection / active-sheet UI state

$wb = $excel.ActiveWorkbook
$wsHeader = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# ---- header sheet: fix row 3 text and highlight it ----
$wsHeader.Range("A3").Value = "smartline: https://services.ga.gov.au/gis/rest/services/Geomorphology_Smartline/MapServer"
$wsHeader.Range("A3").Interior.Color = 65535

# ---- SSSOM sheet: fix the "Soft rock shoes" typo in row 4 ----
$wsData.Range("A4").Value = "smartline:Soft rock shores"
$wsData.Range("B4").Value = "Soft rock shores"

# ---- SSSOM sheet: widen column B to fit the new longer labels ----
$wsData.Columns.Item(2).ColumnWidth = 31

# ---- SSSOM sheet: new rows 7-11 (Rebecca Jordan's additional mappings) ----

# Row 7: Muddy shores -> MT1.2 Muddy Shorelines
$wsData.Range("A7").Value = "smartline:Muddy shores"
$wsData.Range("B7").Value = "Muddy shores"
$wsData.Range("C7").Value = "skos:broadMatch"
$wsData.Range("D7").Value = "get:groups/MT1.2"
$wsData.Range("E7").Value = "MT1.2 Muddy Shorelines"
$wsData.Range("F7").Value = "semapv:ManualMappingCuration"
$wsData.Range("G7").Value = "orcid:0000-0002-4048-6792"
$wsData.Range("H7").Value = "Rebecca Jordan"
$wsData.Range("I7").Value = 45455
$wsData.Range("K7").Value = "status:draft"

# Row 8: Coarse sediment shores -> MT1.4 Boulder and cobble shores
$wsData.Range("A8").Value = "smartline:Coarse sediment shores"
$wsData.Range("B8").Value = "Coarse sediment shores"
$wsData.Range("C8").Value = "skos:broadMatch"
$wsData.Range("D8").Value = "get:groups/MT1.4"
$wsData.Range("E8").Value = "MT1.4 Boulder and cobble shores"
$wsData.Range("F8").Value = "semapv:ManualMappingCuration"
$wsData.Range("G8").Value = "orcid:0000-0002-4048-6792"
$wsData.Range("H8").Value = "Rebecca Jordan"
$wsData.Range("I8").Value = 45455
$wsData.Range("K8").Value = "status:draft"

# Row 9: Undifferentiated sediment shores (no mapping target yet - C/D/E left blank but highlighted)
$wsData.Range("A9").Value = "smartline:Undifferentiated sediment shores"
$wsData.Range("B9").Value = "Undifferentiated sediment shores"
$wsData.Range("C9:E9").Interior.Color = 65535
$wsData.Range("F9").Value = "semapv:ManualMappingCuration"
$wsData.Range("G9").Value = "orcid:0000-0002-4048-6792"
$wsData.Range("H9").Value = "Rebecca Jordan"
$wsData.Range("I9").Value = 45455
$wsData.Range("K9").Value = "status:draft"

# Row 10: Coral coasts (no mapping target yet - C/D/E left blank but highlighted)
$wsData.Range("A10").Value = "smartline:Coral coasts"
$wsData.Range("B10").Value = "Coral coasts"
$wsData.Range("C10:E10").Interior.Color = 65535
$wsData.Range("F10").Value = "semapv:ManualMappingCuration"
$wsData.Range("G10").Value = "orcid:0000-0002-4048-6792"
$wsData.Range("H10").Value = "Rebecca Jordan"
$wsData.Range("I10").Value = 45455
$wsData.Range("K10").Value = "status:draft"

# Row 11: No stability classification -> owl:Nothing / Unclassified
$wsData.Range("A11").Value = "smartline:No stability classification"
$wsData.Range("B11").Value = "No stability classification"
$wsData.Range("C11").Value = "skos:broadMatch"
$wsData.Range("D11").Value = "owl:Nothing"
$wsData.Range("E11").Value = "Unclassified"
$wsData.Range("D11:E11").WrapText = $true
$wsData.Range("F11").Value = "semapv:ManualMappingCuration"
$wsData.Range("F11").Font.Color = 2630431
$wsData.Range("G11").Value = "orcid:0000-0002-4048-6792"
$wsData.Range("H11").Value = "Rebecca Jordan"
$wsData.Range("I11").Value = 45455
$wsData.Range("K11").Value = "status:draft"

# ---- selections / active sheet to match the saved UI state ----
$wsHeader.Range("A3").Select()
$wsData.Range("E22").Select()
$wsData.Activate()
